$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Split "Pohlker" into its own run with proofErr spellStart/spellEnd ---
$pohlkerCell = $t.Cell(4, 3)
$pohlkerXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:pPr>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
<w:kern w:val="0"/>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
<w:lang w:eastAsia="en-GB"/>
<w14:ligatures w14:val="none"/>
</w:rPr>
</w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
<w:kern w:val="0"/>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
<w:lang w:eastAsia="en-GB"/>
<w14:ligatures w14:val="none"/>
</w:rPr>
<w:t>Pohlker</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
<w:kern w:val="0"/>
<w:sz w:val="22"/>
<w:szCs w:val="22"/>
<w:lang w:eastAsia="en-GB"/>
<w14:ligatures w14:val="none"/>
</w:rPr>
<w:t>, 2021, "Respiratory aerosols and droplets in the transmission of infectious diseases"</w:t>
</w:r>
</w:p>
'@
$pohlkerCell.Range.InsertXML($pohlkerXml)

# --- 2. Add new row: Topcu ---
$t.Rows.Add() | Out-Null
$r = $t.Rows.Count
$t.Cell($r, 1).Range.Text = "Introduction to Airborne Disease Transmission Indoors"
$t.Cell($r, 2).Range.Text = "Importance and Impact"
$t.Cell($r, 3).Range.Text = 'Topcu, 2020, "The impact of COVID-19 on emerging stock markets"'

# --- 3. Add new row: Dubey ---
$t.Rows.Add() | Out-Null
$r = $t.Rows.Count
$t.Cell($r, 1).Range.Text = "Introduction to Airborne Disease Transmission Indoors"
$t.Cell($r, 2).Range.Text = "Importance and Impact"
$t.Cell($r, 3).Range.Text = 'Dubey, 2020, "Psychosocial impact of COVID-19"'

Write-Host "Final row count: " $t.Rows.Count
